# Update cryptos list: refresh Price (column D) and Volume(1h) (column E)
# values for several rows, per the latest GitHub Actions scrape.
# Some new Price values look like plain numbers (e.g. "302.01"); since
# this column holds text-formatted price strings (others contain two
# decimal groups like "43.111.32" which cannot be numbers), those cells
# are explicitly set to Text format first so Excel does not auto-convert
# them into floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.111.32"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "2.314.03"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.01"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.03"
$ws.Range("E6").Value = "  -1.94%  "
$ws.Range("E7").Value = "  +1.00%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +1.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.80"
$ws.Range("E10").Value = "  +1.78%  "
$ws.Range("E11").Value = "  -0.76%  "
$ws.Range("E12").Value = "  -1.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.05"
$ws.Range("E13").Value = "  +0.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.91"
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("D15").Value = "2.673.47"
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("D16").Value = "2.324.91"
$ws.Range("E16").Value = "  +0.96%  "
$ws.Range("E17").Value = "  -2.88%  "
$ws.Range("D18").Value = "43.017.85"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.53"
$ws.Range("E19").Value = "  +6.82%  "
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.17"
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.98"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.44"
$ws.Range("E24").Value = "  -1.80%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.44"
$ws.Range("E26").Value = "  -1.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.02"
$ws.Range("E27").Value = "  +1.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.87"
$ws.Range("E28").Value = "  +0.47%  "
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("E30").Value = "  -10.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.61"
$ws.Range("E31").Value = "  -1.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.25"
$ws.Range("E32").Value = "  +4.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.92"
$ws.Range("E33").Value = "  +6.60%  "
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.38"
$ws.Range("E35").Value = "  +7.61%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("E39").Value = "  +0.83%  "
$ws.Range("E40").Value = "  -2.18%  "
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("D42").Value = "1.999.70"
$ws.Range("E42").Value = "  -0.18%  "
$ws.Range("E43").Value = "  +0.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.17"
$ws.Range("E44").Value = "  -4.28%  "
$ws.Range("E45").Value = "  -1.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.56"
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.93"
$ws.Range("E48").Value = "  -0.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "75.04"
$ws.Range("E49").Value = "  +7.07%  "
$ws.Range("D50").Value = "2.539.21"
$ws.Range("E50").Value = "  +0.70%  "
$ws.Range("E51").Value = "  +0.36%  "